$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap data (columns B:AC) between row 26 and row 27
$ws.Range("B26").Value = 6732773
$ws.Range("B27").Value = 6732711
$ws.Range("C26").Value = 'Lithuania A Lyga'
$ws.Range("C27").Value = 'Lithuania A Lyga'
$ws.Range("D26").Value = 'Lithuania A Lyga'
$ws.Range("D27").Value = 'Lithuania A Lyga'
$ws.Range("E26").Value = 45109.58333333334
$ws.Range("E27").Value = 45109.58333333334
$ws.Range("F26").Value = 'Suduva Marijampole'
$ws.Range("F27").Value = 'Banga Gargzdai'
$ws.Range("G26").Value = 'Hegelmann Litauen'
$ws.Range("G27").Value = 'FK Zalgiris Vilnius'
$ws.Range("H26").Value = 0
$ws.Range("H27").Value = 1
$ws.Range("I26").Value = 1
$ws.Range("I27").Value = 4
$ws.Range("J26").Value = 'A'
$ws.Range("J27").Value = 'A'
$ws.Range("K26").Value = 5
$ws.Range("K27").Value = 5
$ws.Range("L26").Value = 3.8
$ws.Range("L27").Value = 3.6
$ws.Range("M26").Value = 1.533
$ws.Range("M27").Value = 1.571
$ws.Range("N26").Value = 5
$ws.Range("N27").Value = 11
$ws.Range("O26").Value = 4.2
$ws.Range("O27").Value = 4.75
$ws.Range("P26").Value = 1.533
$ws.Range("P27").Value = 1.25
$ws.Range("Q26").Value = 1
$ws.Range("Q27").Value = 1.5
$ws.Range("R26").Value = 1.875
$ws.Range("R27").Value = 1.975
$ws.Range("S26").Value = 1.925
$ws.Range("S27").Value = 1.825
$ws.Range("T26").Value = 2.5
$ws.Range("T27").Value = 2.5
$ws.Range("U26").Value = 1.9
$ws.Range("U27").Value = 1.8
$ws.Range("V26").Value = 1.9
$ws.Range("V27").Value = 2
$ws.Range("W26").Value = -1
$ws.Range("W27").Value = -1
$ws.Range("X26").Value = -1
$ws.Range("X27").Value = -1
$ws.Range("Y26").Value = 0.5329999999999999
$ws.Range("Y27").Value = 0.25
$ws.Range("Z26").Value = 0
$ws.Range("Z27").Value = -1
$ws.Range("AA26").Value = 0
$ws.Range("AA27").Value = 0.825
$ws.Range("AB26").Value = -1
$ws.Range("AB27").Value = 0.8
$ws.Range("AC26").Value = 0.8999999999999999
$ws.Range("AC27").Value = -1

# Swap data (columns B:AC) between row 89 and row 90
$ws.Range("B89").Value = 7326568
$ws.Range("B90").Value = 6732827
$ws.Range("C89").Value = 'Lithuania A Lyga'
$ws.Range("C90").Value = 'Lithuania A Lyga'
$ws.Range("D89").Value = 'Lithuania A Lyga'
$ws.Range("D90").Value = 'Lithuania A Lyga'
$ws.Range("E89").Value = 45220.375
$ws.Range("E90").Value = 45220.375
$ws.Range("F89").Value = 'Hegelmann Litauen'
$ws.Range("F90").Value = 'FK Dziugas Telsiai'
$ws.Range("G89").Value = 'Panevezys'
$ws.Range("G90").Value = 'FK Kauno Zalgiris'
$ws.Range("H89").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("I90").Value = 2
$ws.Range("J89").Value = 'D'
$ws.Range("J90").Value = 'A'
$ws.Range("K89").Value = 2.375
$ws.Range("K90").Value = 6
$ws.Range("L89").Value = 3.2
$ws.Range("L90").Value = 3.9
$ws.Range("M89").Value = 2.625
$ws.Range("M90").Value = 1.444
$ws.Range("N89").Value = 2.7
$ws.Range("N90").Value = 4.75
$ws.Range("O89").Value = 3.2
$ws.Range("O90").Value = 3.6
$ws.Range("P89").Value = 2.3
$ws.Range("P90").Value = 1.65
$ws.Range("Q89").Value = 0
$ws.Range("Q90").Value = 0.75
$ws.Range("R89").Value = 2.05
$ws.Range("R90").Value = 1.9
$ws.Range("S89").Value = 1.75
$ws.Range("S90").Value = 1.9
$ws.Range("T89").Value = 2.25
$ws.Range("T90").Value = 2.5
$ws.Range("U89").Value = 1.875
$ws.Range("U90").Value = 1.95
$ws.Range("V89").Value = 1.925
$ws.Range("V90").Value = 1.85
$ws.Range("W89").Value = -1
$ws.Range("W90").Value = -1
$ws.Range("X89").Value = 2.2
$ws.Range("X90").Value = -1
$ws.Range("Y89").Value = -1
$ws.Range("Y90").Value = 0.6499999999999999
$ws.Range("Z89").Value = 0
$ws.Range("Z90").Value = -1
$ws.Range("AA89").Value = 0
$ws.Range("AA90").Value = 0.8999999999999999
$ws.Range("AB89").Value = -1
$ws.Range("AB90").Value = -1
$ws.Range("AC89").Value = 0.925
$ws.Range("AC90").Value = 0.8500000000000001

# Swap data (columns B:AC) between row 101 and row 103
$ws.Range("B101").Value = 6732837
$ws.Range("B103").Value = 6732727
$ws.Range("C101").Value = 'Lithuania A Lyga'
$ws.Range("C103").Value = 'Lithuania A Lyga'
$ws.Range("D101").Value = 'Lithuania A Lyga'
$ws.Range("D103").Value = 'Lithuania A Lyga'
$ws.Range("E101").Value = 45242.41319444445
$ws.Range("E103").Value = 45242.41319444445
$ws.Range("F101").Value = 'Suduva Marijampole'
$ws.Range("F103").Value = 'FK Zalgiris Vilnius'
$ws.Range("G101").Value = 'FK Riteriai'
$ws.Range("G103").Value = 'FK Dainava Alytus'
$ws.Range("H101").Value = 0
$ws.Range("H103").Value = 1
$ws.Range("I101").Value = 3
$ws.Range("I103").Value = 0
$ws.Range("J101").Value = 'A'
$ws.Range("J103").Value = 'H'
$ws.Range("K101").Value = 3.6
$ws.Range("K103").Value = 1.285
$ws.Range("L101").Value = 3.6
$ws.Range("L103").Value = 5.5
$ws.Range("M101").Value = 1.8
$ws.Range("M103").Value = 6.5
$ws.Range("N101").Value = 3
$ws.Range("N103").Value = 1.3
$ws.Range("O101").Value = 3.6
$ws.Range("O103").Value = 5.5
$ws.Range("P101").Value = 2
$ws.Range("P103").Value = 6
$ws.Range("Q101").Value = 0.25
$ws.Range("Q103").Value = -1.5
$ws.Range("R101").Value = 2
$ws.Range("R103").Value = 1.9
$ws.Range("S101").Value = 1.8
$ws.Range("S103").Value = 1.9
$ws.Range("T101").Value = 2.5
$ws.Range("T103").Value = 2.75
$ws.Range("U101").Value = 1.975
$ws.Range("U103").Value = 1.8
$ws.Range("V101").Value = 1.825
$ws.Range("V103").Value = 2
$ws.Range("W101").Value = -1
$ws.Range("W103").Value = 0.3
$ws.Range("X101").Value = -1
$ws.Range("X103").Value = -1
$ws.Range("Y101").Value = 1
$ws.Range("Y103").Value = -1
$ws.Range("Z101").Value = -1
$ws.Range("Z103").Value = -1
$ws.Range("AA101").Value = 0.8
$ws.Range("AA103").Value = 0.8999999999999999
$ws.Range("AB101").Value = 0.9750000000000001
$ws.Range("AB103").Value = -1
$ws.Range("AC101").Value = -1
$ws.Range("AC103").Value = 1

# Swap data (columns B:AC) between row 102 and row 104
$ws.Range("B102").Value = 6732836
$ws.Range("B104").Value = 6732834
$ws.Range("C102").Value = 'Lithuania A Lyga'
$ws.Range("C104").Value = 'Lithuania A Lyga'
$ws.Range("D102").Value = 'Lithuania A Lyga'
$ws.Range("D104").Value = 'Lithuania A Lyga'
$ws.Range("E102").Value = 45242.41319444445
$ws.Range("E104").Value = 45242.41319444445
$ws.Range("F102").Value = 'FK Siauliai'
$ws.Range("F104").Value = 'Panevezys'
$ws.Range("G102").Value = 'Banga Gargzdai'
$ws.Range("G104").Value = 'FK Dziugas Telsiai'
$ws.Range("H102").Value = 3
$ws.Range("H104").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J102").Value = 'H'
$ws.Range("J104").Value = 'D'
$ws.Range("K102").Value = 1.222
$ws.Range("K104").Value = 1.25
$ws.Range("L102").Value = 5.5
$ws.Range("L104").Value = 5.5
$ws.Range("M102").Value = 9
$ws.Range("M104").Value = 7.5
$ws.Range("N102").Value = 1.363
$ws.Range("N104").Value = 1.45
$ws.Range("O102").Value = 4.5
$ws.Range("O104").Value = 4.5
$ws.Range("P102").Value = 7
$ws.Range("P104").Value = 5
$ws.Range("Q102").Value = -1.25
$ws.Range("Q104").Value = -1
$ws.Range("R102").Value = 1.9
$ws.Range("R104").Value = 1.775
$ws.Range("S102").Value = 1.9
$ws.Range("S104").Value = 2.025
$ws.Range("T102").Value = 2.5
$ws.Range("T104").Value = 2.5
$ws.Range("U102").Value = 1.975
$ws.Range("U104").Value = 1.875
$ws.Range("V102").Value = 1.825
$ws.Range("V104").Value = 1.925
$ws.Range("W102").Value = 0.363
$ws.Range("W104").Value = -1
$ws.Range("X102").Value = -1
$ws.Range("X104").Value = 3.5
$ws.Range("Y102").Value = -1
$ws.Range("Y104").Value = -1
$ws.Range("Z102").Value = 0.8999999999999999
$ws.Range("Z104").Value = -1
$ws.Range("AA102").Value = -1
$ws.Range("AA104").Value = 1.025
$ws.Range("AB102").Value = 0.9750000000000001
$ws.Range("AB104").Value = -1
$ws.Range("AC102").Value = -1
$ws.Range("AC104").Value = 0.925

